$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.945012333333333
$ws.Range("H2").Value = 5.835037
$ws.Range("Q2").Value = 1.063705201626889
$ws.Range("R2").Value = 9.573346814641999
